$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new blank column before column N ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N").Insert()
$wsRepay.Columns("N").ColumnWidth = $wsRepay.Columns("M").ColumnWidth

# --- Acc_Upfront sheet: move selection, no longer the active/selected tab ---
$wsUpfront = $wb.Worksheets.Item("Acc_Upfront")
$wsUpfront.Range("E7").Select() | Out-Null

# --- Make "Repayment schedule" the active sheet/tab with new selection ---
$wsRepay.Activate() | Out-Null
$wsRepay.Range("S7").Select() | Out-Null
